$d = $word.ActiveDocument

$replacements = @(
    @("2025-10-11 Saturday", "2025-10-12 Sunday"),
    @("371×3=", "169×6="),
    @("107×4=", "228×7="),
    @("751×8=", "162×2="),
    @("514×5=", "264×3="),
    @("561×8=", "319×7="),
    @("500×2=", "688×4="),
    @("995×2=", "786×9="),
    @("778×6=", "794×9="),
    @("810×7=", "486×2="),
    @("963×6=", "875×2="),
    @("261×9=", "743×2="),
    @("958×4=", "214×6="),
    @("127×9=", "811×8="),
    @("773×3=", "824×5="),
    @("747×3=", "158×7="),
    @("351×9=", "164×2="),
    @("139×8=", "276×8="),
    @("520×9=", "779×3="),
    @("894×6=", "728×5="),
    @("873×7=", "358×2="),
    @("566×3=", "397×5="),
    @("536×2=", "162×9="),
    @("465×7=", "106×7="),
    @("773×8=", "367×6="),
    @("376×7=", "711×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
